# Refresh cached Market Board price/profit figures across all Leve sheets
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) per the latest
# scheduled market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 645.8276
$ws.Range("I28").Value = 686.6
$ws.Range("J28").Value = 391
$ws.Range("K28").Value = 686.6
$ws.Range("L28").Value = 391
$ws.Range("M28").Value = -201.6
$ws.Range("N28").Value = -1361
$ws.Range("H40").Value = 1347.3077
$ws.Range("I40").Value = 1187.5862
$ws.Range("J40").Value = 1548.6957
$ws.Range("K40").Value = 1187.5862
$ws.Range("L40").Value = 1548.6957
$ws.Range("M40").Value = -1012.5862
$ws.Range("N40").Value = -1898.6957
$ws.Range("H62").Value = 13769.125
$ws.Range("I62").Value = 1450.4286
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 1450.4286
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -826.4286
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 13769.125
$ws.Range("I65").Value = 1450.4286
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 7252.143
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -4132.143
$ws.Range("N65").Value = -506240
$ws.Range("H96").Value = 722
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 722
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2166
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4912
$ws.Range("H100").Value = 20001684
$ws.Range("I100").Value = 1461.1177
$ws.Range("J100").Value = 62502156
$ws.Range("K100").Value = 1461.1177
$ws.Range("L100").Value = 62502156
$ws.Range("M100").Value = -920.1177
$ws.Range("N100").Value = -62503238
$ws.Range("H103").Value = 8000470.5
$ws.Range("I103").Value = 330.7857
$ws.Range("J103").Value = 18182466
$ws.Range("K103").Value = 992.3571000000001
$ws.Range("L103").Value = 54547398
$ws.Range("M103").Value = -406.3571000000001
$ws.Range("N103").Value = -54548570
$ws.Range("H107").Value = 747.09375
$ws.Range("I107").Value = 849.11536
$ws.Range("J107").Value = 305
$ws.Range("K107").Value = 849.11536
$ws.Range("L107").Value = 305
$ws.Range("M107").Value = 1070.88464
$ws.Range("N107").Value = -4145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2692.353
$ws.Range("I2").Value = 1666.9231
$ws.Range("J2").Value = 6025
$ws.Range("K2").Value = 1666.9231
$ws.Range("L2").Value = 6025
$ws.Range("M2").Value = -1553.9231
$ws.Range("N2").Value = -6251
$ws.Range("H45").Value = 4843.875
$ws.Range("I45").Value = 4157.2856
$ws.Range("J45").Value = 5377.8887
$ws.Range("K45").Value = 4157.2856
$ws.Range("L45").Value = 5377.8887
$ws.Range("M45").Value = -3780.2856
$ws.Range("N45").Value = -6131.8887
$ws.Range("H63").Value = 2424.7083
$ws.Range("I63").Value = 1478.0714
$ws.Range("J63").Value = 3750
$ws.Range("K63").Value = 1478.0714
$ws.Range("L63").Value = 3750
$ws.Range("M63").Value = -792.0714
$ws.Range("N63").Value = -5122
$ws.Range("H66").Value = 2424.7083
$ws.Range("I66").Value = 1478.0714
$ws.Range("J66").Value = 3750
$ws.Range("K66").Value = 7390.357
$ws.Range("L66").Value = 18750
$ws.Range("M66").Value = -3958.357
$ws.Range("N66").Value = -25614
$ws.Range("H97").Value = 1286.0605
$ws.Range("I97").Value = 1099.4546
$ws.Range("J97").Value = 1659.2727
$ws.Range("K97").Value = 1099.4546
$ws.Range("L97").Value = 1659.2727
$ws.Range("M97").Value = -603.4546
$ws.Range("N97").Value = -2651.2727
$ws.Range("H112").Value = 9392.333000000001
$ws.Range("J112").Value = 9392.333000000001
$ws.Range("L112").Value = 9392.333000000001
$ws.Range("N112").Value = -12346.333
$ws.Range("H114").Value = 35000
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678
$ws.Range("H116").Value = 2692.353
$ws.Range("I116").Value = 1666.9231
$ws.Range("J116").Value = 6025
$ws.Range("K116").Value = 1666.9231
$ws.Range("L116").Value = 6025
$ws.Range("M116").Value = 627.0769
$ws.Range("N116").Value = -10613
$ws.Range("H119").Value = 20475.143
$ws.Range("J119").Value = 20475.143
$ws.Range("L119").Value = 20475.143
$ws.Range("N119").Value = -30151.143
$ws.Range("H121").Value = 15709.615
$ws.Range("J121").Value = 15709.615
$ws.Range("L121").Value = 15709.615
$ws.Range("N121").Value = -19203.615
$ws.Range("H122").Value = 2048.2942
$ws.Range("I122").Value = 1965.5385
$ws.Range("K122").Value = 5896.6155
$ws.Range("M122").Value = -3446.6155
$ws.Range("H132").Value = 1936.4706
$ws.Range("I132").Value = 1684
$ws.Range("J132").Value = 2910.2856
$ws.Range("K132").Value = 5052
$ws.Range("L132").Value = 8730.856800000001
$ws.Range("M132").Value = -2522
$ws.Range("N132").Value = -13790.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2692.353
$ws.Range("I3").Value = 1666.9231
$ws.Range("J3").Value = 6025
$ws.Range("K3").Value = 1666.9231
$ws.Range("L3").Value = 6025
$ws.Range("M3").Value = -1552.9231
$ws.Range("N3").Value = -6253
$ws.Range("H134").Value = 1426.1923
$ws.Range("I134").Value = 1277.85
$ws.Range("J134").Value = 1920.6666
$ws.Range("K134").Value = 3833.55
$ws.Range("L134").Value = 5761.9998
$ws.Range("M134").Value = -1298.55
$ws.Range("N134").Value = -10831.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 980
$ws.Range("I23").Value = 980
$ws.Range("K23").Value = 980
$ws.Range("M23").Value = -740
$ws.Range("H27").Value = 980
$ws.Range("I27").Value = 980
$ws.Range("K27").Value = 980
$ws.Range("M27").Value = -788
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 5666.6665
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 16999.9995
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -21899.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4408.619
$ws.Range("I109").Value = 2996.2307
$ws.Range("J109").Value = 6703.75
$ws.Range("K109").Value = 8988.6921
$ws.Range("L109").Value = 20111.25
$ws.Range("M109").Value = -7948.6921
$ws.Range("N109").Value = -22191.25
$ws.Range("H117").Value = 1352.5714
$ws.Range("I117").Value = 893.4
$ws.Range("J117").Value = 2500.5
$ws.Range("K117").Value = 2680.2
$ws.Range("L117").Value = 7501.5
$ws.Range("M117").Value = 761.8000000000002
$ws.Range("N117").Value = -14385.5
$ws.Range("H123").Value = 5461.857
$ws.Range("J123").Value = 5461.857
$ws.Range("L123").Value = 16385.571
$ws.Range("N123").Value = -21285.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1811.3684
$ws.Range("I97").Value = 1998.9286
$ws.Range("J97").Value = 1286.2
$ws.Range("K97").Value = 1998.9286
$ws.Range("L97").Value = 1286.2
$ws.Range("M97").Value = -1502.9286
$ws.Range("N97").Value = -2278.2
$ws.Range("H113").Value = 1545.1538
$ws.Range("I113").Value = 1458.9375
$ws.Range("J113").Value = 1683.1
$ws.Range("K113").Value = 1458.9375
$ws.Range("L113").Value = 1683.1
$ws.Range("M113").Value = 711.0625
$ws.Range("N113").Value = -6023.1
$ws.Range("H122").Value = 1905.8
$ws.Range("I122").Value = 1852.8462
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 5558.5386
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3108.5386
$ws.Range("N122").Value = -11650
$ws.Range("H126").Value = 11426.931
$ws.Range("I126").Value = 3027.3635
$ws.Range("J126").Value = 16560
$ws.Range("K126").Value = 9082.0905
$ws.Range("L126").Value = 49680
$ws.Range("M126").Value = -6612.0905
$ws.Range("N126").Value = -54620
$ws.Range("H132").Value = 2673.9688
$ws.Range("I132").Value = 2093.8235
$ws.Range("J132").Value = 3331.4666
$ws.Range("K132").Value = 6281.470499999999
$ws.Range("L132").Value = 9994.399800000001
$ws.Range("M132").Value = -3751.470499999999
$ws.Range("N132").Value = -15054.3998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2462
$ws.Range("I61").Value = 1785.4
$ws.Range("J61").Value = 3428.5715
$ws.Range("K61").Value = 1785.4
$ws.Range("L61").Value = 3428.5715
$ws.Range("M61").Value = -1583.4
$ws.Range("N61").Value = -3832.5715
$ws.Range("H113").Value = 2462
$ws.Range("I113").Value = 1785.4
$ws.Range("J113").Value = 3428.5715
$ws.Range("K113").Value = 1785.4
$ws.Range("L113").Value = 3428.5715
$ws.Range("M113").Value = 384.5999999999999
$ws.Range("N113").Value = -7768.5715
$ws.Range("H122").Value = 2830.762
$ws.Range("I122").Value = 2118.1
$ws.Range("J122").Value = 3478.6365
$ws.Range("K122").Value = 6354.299999999999
$ws.Range("L122").Value = 10435.9095
$ws.Range("M122").Value = -3904.299999999999
$ws.Range("N122").Value = -15335.9095
$ws.Range("H132").Value = 3547.5862
$ws.Range("I132").Value = 3262.818
$ws.Range("K132").Value = 9788.454000000002
$ws.Range("M132").Value = -7258.454000000002
$ws.Range("H136").Value = 2539.52
$ws.Range("I136").Value = 1827.1111
$ws.Range("J136").Value = 4371.4287
$ws.Range("K136").Value = 5481.3333
$ws.Range("L136").Value = 13114.2861
$ws.Range("M136").Value = -2931.3333
$ws.Range("N136").Value = -18214.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 42778
$ws.Range("J105").Value = 42778
$ws.Range("L105").Value = 42778
$ws.Range("N105").Value = -49766
$ws.Range("H107").Value = 1470.3448
$ws.Range("I107").Value = 938.8421
$ws.Range("J107").Value = 2480.2
$ws.Range("K107").Value = 2816.5263
$ws.Range("L107").Value = 7440.599999999999
$ws.Range("M107").Value = -896.5263
$ws.Range("N107").Value = -11280.6
$ws.Range("H113").Value = 127017.125
$ws.Range("I113").Value = 200467.4
$ws.Range("K113").Value = 601402.2
$ws.Range("M113").Value = -599232.2
$ws.Range("H122").Value = 2528.7917
$ws.Range("I122").Value = 2117.647
$ws.Range("J122").Value = 3527.2856
$ws.Range("K122").Value = 6352.941
$ws.Range("L122").Value = 10581.8568
$ws.Range("M122").Value = -3902.941
$ws.Range("N122").Value = -15481.8568
$ws.Range("H138").Value = 48206.5
$ws.Range("J138").Value = 48206.5
$ws.Range("L138").Value = 48206.5
$ws.Range("N138").Value = -58486.5

Write-Output "Applied 265 cell updates across 8 sheets."
